$wb = $excel.ActiveWorkbook

# Updated "want to go" counts (column F) for the rows that remain after the
# first data row (the 2024.01.13 "秋绥冬禧国乙only" entry) is removed.
$newCounts = @{
    2 = 1905
    3 = 581
    4 = 1257
    5 = 6339
    6 = 154
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the obsolete first data row (2024.01.13, row 2). Remaining rows
    # shift up by one automatically.
    $ws.Rows.Item(2).Delete()

    # Renumber the leading index column (A) for the remaining data rows.
    for ($r = 2; $r -le 6; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Apply the updated "想去人数" (column F) values.
    foreach ($r in $newCounts.Keys) {
        $ws.Cells.Item($r, 6).Value = $newCounts[$r]
    }
}
